# Refresh the cryptos price/volume snapshot (GitHub Actions data pull).
# Numeric-looking "Price" values are entered with a leading apostrophe so
# Excel keeps them as literal text (matching the sheet's existing
# plain-text D-column cells, e.g. "13.00" instead of being normalized to 13).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.468.86'
$ws.Range('E2').Value = '  +2.02%  '
$ws.Range('D3').Value = '3.153.95'
$ws.Range('E3').Value = '  +2.57%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''535.96'
$ws.Range('E5').Value = '  +2.11%  '
$ws.Range('D6').Value = '''139.87'
$ws.Range('E6').Value = '  +2.87%  '
$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '''0.512'
$ws.Range('E8').Value = '  +8.93%  '
$ws.Range('D9').Value = '''7.33'
$ws.Range('E9').Value = '  +1.40%  '
$ws.Range('E10').Value = '  +3.41%  '
$ws.Range('E11').Value = '  +4.99%  '
$ws.Range('E12').Value = '  +2.00%  '
$ws.Range('D13').Value = '3.695.50'
$ws.Range('E13').Value = '  +2.55%  '
$ws.Range('D14').Value = '''25.74'
$ws.Range('E14').Value = '  +2.52%  '
$ws.Range('E15').Value = '  +6.45%  '
$ws.Range('D16').Value = '58.522.75'
$ws.Range('E16').Value = '  +2.01%  '
$ws.Range('D17').Value = '3.164.92'
$ws.Range('E17').Value = '  +2.57%  '
$ws.Range('D18').Value = '''6.22'
$ws.Range('E18').Value = '  +6.55%  '
$ws.Range('D19').Value = '''13.00'
$ws.Range('E19').Value = '  +4.61%  '
$ws.Range('D20').Value = '''8.21'
$ws.Range('E20').Value = '  +4.90%  '
$ws.Range('D21').Value = '''371.94'
$ws.Range('E21').Value = '  +6.98%  '
$ws.Range('D22').Value = '''5.78'
$ws.Range('E22').Value = '  +1.44%  '
$ws.Range('D23').Value = '''1.00'
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('D24').Value = '''69.92'
$ws.Range('E24').Value = '  +2.11%  '
$ws.Range('E25').Value = '  +3.33%  '
$ws.Range('E26').Value = '  +0.97%  '
$ws.Range('E27').Value = '  -0.24%  '
$ws.Range('E28').Value = '  +13.84%  '
$ws.Range('D29').Value = '0.0₃0870'
$ws.Range('E29').Value = '  +3.50%  '
$ws.Range('E30').Value = '  +2.45%  '
$ws.Range('E31').Value = '  +5.35%  '
$ws.Range('E32').Value = '  +4.70%  '
$ws.Range('E33').Value = '  +7.97%  '
$ws.Range('E34').Value = '  +4.65%  '
$ws.Range('D35').Value = '''159.72'
$ws.Range('E35').Value = '  +0.78%  '
$ws.Range('D36').Value = '''6.22'
$ws.Range('E36').Value = '  +4.20%  '
$ws.Range('D37').Value = '''1.36'
$ws.Range('E37').Value = '  +11.65%  '
$ws.Range('D38').Value = '''25.25'
$ws.Range('E38').Value = '  -0.57%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Value = '''1.68'
$ws.Range('E39').Value = '  +7.21%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').Value = '2.643.71'
$ws.Range('E40').Value = '  +9.45%  '
$ws.Range('D41').Value = '''0.0682'
$ws.Range('E41').Value = '  +3.97%  '
$ws.Range('D42').Value = '''4.17'
$ws.Range('E42').Value = '  +3.98%  '
$ws.Range('D43').Value = '''38.79'
$ws.Range('D44').Value = '''0.708'
$ws.Range('E44').Value = '  +2.90%  '
$ws.Range('E45').Value = '  +8.55%  '
$ws.Range('D46').Value = '''1.00'
$ws.Range('E46').Value = '  -0.02%  '
$ws.Range('D47').Value = '3.195.80'
$ws.Range('E47').Value = '  +2.64%  '
$ws.Range('B48').Value = 'ONDO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D48').Value = '''0.983'
$ws.Range('E48').Value = '  +4.88%  '
$ws.Range('D49').Value = '''6.21'
$ws.Range('E49').Value = '  +4.06%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').Value = '''0.101'
$ws.Range('E50').Value = '  +10.26%  '
$ws.Range('D51').Value = '''20.24'
$ws.Range('E51').Value = '  +4.76%  '
